$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 and 4 entirely (events no longer present after this update)
$ws.Range("A3:L4").EntireRow.Delete()

# Row 2: replace with the new event record
$ws.Range("A2").Value = '281474990867465-1749055419927'
$ws.Range("B2").Value = 'Harsh Turn'
$ws.Range("C2").Value = '2025-06-04T10:43:39.927'

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '281474990867465'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '124'
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = '52215846'
$ws.Range("F2").Style = "Normal"

$ws.Range("G2").Value = 'ALVARO ZAPATA'

$ws.Range("H2").Value = 20.78436023
$ws.Range("I2").Value = -105.36537465
$ws.Range("J2").Value = 0.484667956829071

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990867465/1749055414927/sdxp5cjHuF-camera-video-segment-driver-1749055419927.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK6N6IWR%2F20250605%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250605T170108Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG4aCXVzLXdlc3QtMiJHMEUCIBwEb2kFMWTALqE0zWBKSdJV6tCQnkIxVPcm3YI0jBjdAiEAmmD6r22uXRqwiWAMueF%2FX%2BYHbRNBqazt5aloM4Wt%2Fhwq3QMIRxAEGgw3ODEyMDQ5NDIyNDQiDG9YZU02CweWYUJFYyq6A%2FHF3lV1oL0OMCdS7hQwGbJNnSEXmUOAqYRQc2lXuC%2BjJmGsnZf2DnB1XogKVtMNiU81A5OJf0DAMhhu38jh2xjKanrq8NiiRycApLMQNm49pAOy%2B6Wd06zlZeHunyzhhOI4ZOGyoFFpQbtMkJATRauvh%2FtO3au7Rz%2FU6f9%2FFC9Hqtgy5gbP75wqzYfQ8%2FKWt4iLu6UOCJSS4TOfoEolFo%2FhqArxxPcPSDtLJsw%2FLxJchXT8kaIOFs3l5%2FsuFBBnEd%2FdSzyb1aAoVOuganibPhbvkP1EaIM%2FtHSkpyj3%2FgauDwdmFNXKM5FZ4yZnu7YxXf458OfNL%2F0DltWtBT0is3BE6c%2B%2Bd8PUztwRDuOoNsAlR6sOzhIGnBn3sI%2F6%2BM4ynEIDb%2BARvnVkcHnFO9EeSRK1%2Fm6L9lBcUPnt9ycUHCnCle9tsQrL%2FD%2FBHBzqDYw7WMBx%2FO1V1F0mv8cMpk5%2BvVQi3Q3LkAAfbAE9mgH9Z4eE0Tt4x2XVwuV0a4YJvMlRL125rfzq8E6YAvLw7UztGlTQDCMlSyYSbRqLdMp7BfPoa5m33GSiIdotzRmr8HdTisoQ%2B6d11m3xUhYwj8mGwgY6pQGkMro56aHXPLRWKRZKaG14w3Niq8urcQT7lbsRMugMos0k0ZvbgoNaJqhSrSfTAFdnush3z2%2FJn%2BP4t8K%2Ba6irGA9lKSs9lMFI767XcL7qlcx1COv64WWwpZ9N3Xv48aRNrXMHkDPRfQ5HyOcgRGw5Q8GJOOG7627he2mCp38CD3tHfejSpwKwF8SNPEyTfsoNBWemd4FeLecaLCsKcVpVuE4UtUs%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2006%20Jun%202025%2001%3A01%3A08%20GMT&X-Amz-Signature=e6da25dbb82177777f7cdb3cf1b7c540170d3047323c08f2211526890819f9c2'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990867465/1749055414927/kyWdDmegaS-camera-video-segment-1749055419927.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSMK6N6IWR%2F20250605%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250605T170108Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEG4aCXVzLXdlc3QtMiJHMEUCIBwEb2kFMWTALqE0zWBKSdJV6tCQnkIxVPcm3YI0jBjdAiEAmmD6r22uXRqwiWAMueF%2FX%2BYHbRNBqazt5aloM4Wt%2Fhwq3QMIRxAEGgw3ODEyMDQ5NDIyNDQiDG9YZU02CweWYUJFYyq6A%2FHF3lV1oL0OMCdS7hQwGbJNnSEXmUOAqYRQc2lXuC%2BjJmGsnZf2DnB1XogKVtMNiU81A5OJf0DAMhhu38jh2xjKanrq8NiiRycApLMQNm49pAOy%2B6Wd06zlZeHunyzhhOI4ZOGyoFFpQbtMkJATRauvh%2FtO3au7Rz%2FU6f9%2FFC9Hqtgy5gbP75wqzYfQ8%2FKWt4iLu6UOCJSS4TOfoEolFo%2FhqArxxPcPSDtLJsw%2FLxJchXT8kaIOFs3l5%2FsuFBBnEd%2FdSzyb1aAoVOuganibPhbvkP1EaIM%2FtHSkpyj3%2FgauDwdmFNXKM5FZ4yZnu7YxXf458OfNL%2F0DltWtBT0is3BE6c%2B%2Bd8PUztwRDuOoNsAlR6sOzhIGnBn3sI%2F6%2BM4ynEIDb%2BARvnVkcHnFO9EeSRK1%2Fm6L9lBcUPnt9ycUHCnCle9tsQrL%2FD%2FBHBzqDYw7WMBx%2FO1V1F0mv8cMpk5%2BvVQi3Q3LkAAfbAE9mgH9Z4eE0Tt4x2XVwuV0a4YJvMlRL125rfzq8E6YAvLw7UztGlTQDCMlSyYSbRqLdMp7BfPoa5m33GSiIdotzRmr8HdTisoQ%2B6d11m3xUhYwj8mGwgY6pQGkMro56aHXPLRWKRZKaG14w3Niq8urcQT7lbsRMugMos0k0ZvbgoNaJqhSrSfTAFdnush3z2%2FJn%2BP4t8K%2Ba6irGA9lKSs9lMFI767XcL7qlcx1COv64WWwpZ9N3Xv48aRNrXMHkDPRfQ5HyOcgRGw5Q8GJOOG7627he2mCp38CD3tHfejSpwKwF8SNPEyTfsoNBWemd4FeLecaLCsKcVpVuE4UtUs%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2006%20Jun%202025%2001%3A01%3A08%20GMT&X-Amz-Signature=c94d005450f6422ef1dd8b496a7b4becdc4927680c749f8dfc0e94fb0fc9ad9e'
